$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R mirrors column Q's formatting for rows 4-6, with new data values.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2022

$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 8.6821914120339212

$ws.Range("Q6").Copy($ws.Range("R6"))
$ws.Range("R6").Value = 12.221423436376707

# Match the author's final selection (S4) recorded in the saved view state.
$ws.Range("S4").Select()
